$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$query = @'
MATCH (samp:sample)
WHERE samp.specific_sample_pathology IN ['T Cell Lymphoma']
MATCH (samp)-->(c:case)-->(s:study)
MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (f)-->(parent)
OPTIONAL MATCH (c)-->(cv:canine_individual)
WITH
	DISTINCT f, samp, c, s, parent, cv,
  	['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
	toInteger(floor(log(f.file_size)/log(1024))) as i,
	2 as precision
 WITH
  	samp, c, s, f, parent, cv,
    f.file_size /(1024^i) AS value,
    10^precision AS factor,
    units[i] as unit
WITH
    samp, c, s, f, unit, parent, cv,
    round(factor * value)/factor AS size
RETURN
	coalesce(f.file_name, '') AS `File Name`,
	coalesce(f.file_format, '') AS `Format`,
	coalesce(f.file_type, '') AS `File Type`,
	CASE size % 1
	  WHEN 0
	  THEN apoc.convert.toInteger(size)+' ' +unit
	  ELSE size+' ' +unit
	END AS Size,
	head(labels(parent)) AS `Association`,
	coalesce(f.file_description,'') AS `Description`,
	coalesce(samp.sample_id, '') AS `Sample ID`,
	coalesce(c.case_id,'') as `Case ID`,
	coalesce(cv.canine_individual_id,'') AS `Canine ID`,
	CASE
	  WHEN s.clinical_study_designation IS NULL
	  THEN parent.clinical_study_designation
	  ELSE s.clinical_study_designation
	END AS `Study Code`
ORDER BY `File Name`
LIMIT 100
'@

$ws.Range("D2").Value = $query
$ws.Range("D3").Value = $query
$ws.Range("D4").Value = $query
$ws.Range("D5").Value = $query
